# Append: 2026-01-06 06:40 JST
# Prepend two new scraped listings to the "ランサーズ" sheet, pushing the
# existing rows down by two and refreshing the "取得日時" timestamp on
# every row to the new run time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$newTimestamp = "2026-01-06 06:40:29"

# --- 1. Make room: insert two fresh rows above the current row 2 -----------
$ws.Range("A2:A3").EntireRow.Insert()

# --- 2. Refresh the timestamp column for every data row (2-7) --------------
$ws.Range("A2").Value = $newTimestamp
$ws.Range("A3").Value = $newTimestamp
$ws.Range("A4").Value = $newTimestamp
$ws.Range("A5").Value = $newTimestamp
$ws.Range("A6").Value = $newTimestamp
$ws.Range("A7").Value = $newTimestamp

# --- 3. Fill in the two brand-new listings (row 2 and row 3) ---------------
$url2 = "https://www.lancers.jp/work/detail/5445159"
$url3 = "https://www.lancers.jp/work/detail/5445154"
$url4 = "https://www.lancers.jp/work/detail/5425629"
$url5 = "https://www.lancers.jp/work/detail/5465878"
$url6 = "https://www.lancers.jp/work/detail/5465685"
$url7 = "https://www.lancers.jp/work/detail/5465836"

$ws.Range("B2").Value = "法人向け生成AIサービス(RAG・議事録機能)の設計・開発を支援エンジニア募集(AI/バックエンド)"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = $url2
$ws.Range("G2").Value = 368
$ws.Range("H2").Value = "🔥AI,Ai ◆開発"

$ws.Range("B3").Value = "B2B向け生成AIサービス(チャット・RAG)の新規開発プロジェクト推進を支援してくださるPM募集"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "200,000 円 ~ 300,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = $url3
$ws.Range("G3").Value = 368
$ws.Range("H3").Value = "🔥AI,Ai ◆開発"

# --- 3b. The shifted-down rows (now 4-7) already carry their original
#         F-column URL text via the row insert; restate them explicitly so
#         the script is self-contained regardless of insert/shift quirks. --
$ws.Range("F4").Value = $url4
$ws.Range("F5").Value = $url5
$ws.Range("F6").Value = $url6
$ws.Range("F7").Value = $url7

# --- 4. Rebuild the hyperlinks on column F (inserting rows does not carry
#        the hyperlink collection along with the shifted cells) ------------
$ws.Hyperlinks.Delete()

$ws.Hyperlinks.Add($ws.Range("F2"), $url2)
$ws.Range("F2").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F3"), $url3)
$ws.Range("F3").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F4"), $url4)
$ws.Range("F4").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F5"), $url5)
$ws.Range("F5").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F6"), $url6)
$ws.Range("F6").Style = "Hyperlink"

$ws.Hyperlinks.Add($ws.Range("F7"), $url7)
$ws.Range("F7").Style = "Hyperlink"
